$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Age 0-4")
$ws2 = $wb.Worksheets.Item("Age 5-17")
$ws3 = $wb.Worksheets.Item("Age 18-65")
$ws4 = $wb.Worksheets.Item("Age 66+")

# ----- Age 0-4 -----
$ws1.Cells.Item(3, 1).Value = "ALEXANDER"
$ws1.Cells.Item(3, 2).Value = 45
$ws1.Cells.Item(3, 3).Value = 53
$ws1.Cells.Item(3, 4).Value = 0
$ws1.Cells.Item(3, 5).Value = 1
$ws1.Cells.Item(3, 6).Value = 1
$ws1.Cells.Item(3, 7).Value = 2
$ws1.Cells.Item(3, 8).Value = 54

$ws1.Cells.Item(4, 1).Value = "BOND"
$ws1.Cells.Item(4, 2).Value = 90
$ws1.Cells.Item(4, 3).Value = 11
$ws1.Cells.Item(4, 4).Value = 0
$ws1.Cells.Item(4, 5).Value = 0
$ws1.Cells.Item(4, 6).Value = 0
$ws1.Cells.Item(4, 7).Value = 2
$ws1.Cells.Item(4, 8).Value = 56

$ws1.Cells.Item(5, 1).Value = "BOONE"
$ws1.Cells.Item(5, 2).Value = 371
$ws1.Cells.Item(5, 3).Value = 49
$ws1.Cells.Item(5, 4).Value = 3
$ws1.Cells.Item(5, 5).Value = 2
$ws1.Cells.Item(5, 6).Value = 0
$ws1.Cells.Item(5, 7).Value = 188
$ws1.Cells.Item(5, 8).Value = 441

$ws1.Cells.Item(6, 1).Value = "BROWN"
$ws1.Cells.Item(6, 2).Value = 41
$ws1.Cells.Item(6, 3).Value = 0
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(6, 5).Value = 0
$ws1.Cells.Item(6, 6).Value = 0
$ws1.Cells.Item(6, 7).Value = 0
$ws1.Cells.Item(6, 8).Value = 12

$ws1.Cells.Item(7, 1).Value = "BUREAU"
$ws1.Cells.Item(7, 2).Value = 297
$ws1.Cells.Item(7, 3).Value = 12
$ws1.Cells.Item(7, 4).Value = 0
$ws1.Cells.Item(7, 5).Value = 2
$ws1.Cells.Item(7, 6).Value = 0
$ws1.Cells.Item(7, 7).Value = 68
$ws1.Cells.Item(7, 8).Value = 213

$ws1.Cells.Item(104, 1).Value = "Sum:"
$ws1.Cells.Item(104, 2).Value = 66186
$ws1.Cells.Item(104, 3).Value = 51936
$ws1.Cells.Item(104, 4).Value = 557
$ws1.Cells.Item(104, 5).Value = 2936
$ws1.Cells.Item(104, 6).Value = 415
$ws1.Cells.Item(104, 7).Value = 32471
$ws1.Cells.Item(104, 8).Value = 91836

# ----- Age 5-17 -----
$ws2.Cells.Item(3, 1).Value = "ALEXANDER"
$ws2.Cells.Item(3, 2).Value = 127
$ws2.Cells.Item(3, 3).Value = 117
$ws2.Cells.Item(3, 4).Value = 0
$ws2.Cells.Item(3, 5).Value = 0
$ws2.Cells.Item(3, 6).Value = 0
$ws2.Cells.Item(3, 7).Value = 3
$ws2.Cells.Item(3, 8).Value = 198

$ws2.Cells.Item(4, 1).Value = "BOND"
$ws2.Cells.Item(4, 2).Value = 320
$ws2.Cells.Item(4, 3).Value = 36
$ws2.Cells.Item(4, 4).Value = 1
$ws2.Cells.Item(4, 5).Value = 2
$ws2.Cells.Item(4, 6).Value = 0
$ws2.Cells.Item(4, 7).Value = 4
$ws2.Cells.Item(4, 8).Value = 178

$ws2.Cells.Item(5, 1).Value = "BOONE"
$ws2.Cells.Item(5, 2).Value = 796
$ws2.Cells.Item(5, 3).Value = 81
$ws2.Cells.Item(5, 4).Value = 5
$ws2.Cells.Item(5, 5).Value = 6
$ws2.Cells.Item(5, 6).Value = 0
$ws2.Cells.Item(5, 7).Value = 418
$ws2.Cells.Item(5, 8).Value = 1225

$ws2.Cells.Item(6, 1).Value = "BROWN"
$ws2.Cells.Item(6, 2).Value = 109
$ws2.Cells.Item(6, 3).Value = 2
$ws2.Cells.Item(6, 4).Value = 0
$ws2.Cells.Item(6, 5).Value = 0
$ws2.Cells.Item(6, 6).Value = 0
$ws2.Cells.Item(6, 7).Value = 2
$ws2.Cells.Item(6, 8).Value = 29

$ws2.Cells.Item(7, 1).Value = "BUREAU"
$ws2.Cells.Item(7, 2).Value = 593
$ws2.Cells.Item(7, 3).Value = 31
$ws2.Cells.Item(7, 4).Value = 3
$ws2.Cells.Item(7, 5).Value = 8
$ws2.Cells.Item(7, 6).Value = 3
$ws2.Cells.Item(7, 7).Value = 126
$ws2.Cells.Item(7, 8).Value = 546

$ws2.Cells.Item(104, 1).Value = "Sum:"
$ws2.Cells.Item(104, 2).Value = 137609
$ws2.Cells.Item(104, 3).Value = 88743
$ws2.Cells.Item(104, 4).Value = 1144
$ws2.Cells.Item(104, 5).Value = 6029
$ws2.Cells.Item(104, 6).Value = 590
$ws2.Cells.Item(104, 7).Value = 73029
$ws2.Cells.Item(104, 8).Value = 259298

# ----- Age 18-65 -----
$ws3.Cells.Item(9, 1).Value = "ALEXANDER"
$ws3.Cells.Item(9, 2).Value = 473
$ws3.Cells.Item(9, 3).Value = 475
$ws3.Cells.Item(9, 4).Value = 3
$ws3.Cells.Item(9, 5).Value = 1
$ws3.Cells.Item(9, 6).Value = 1
$ws3.Cells.Item(9, 7).Value = 7
$ws3.Cells.Item(9, 8).Value = 55

$ws3.Cells.Item(10, 1).Value = "BOND"
$ws3.Cells.Item(10, 2).Value = 936
$ws3.Cells.Item(10, 3).Value = 74
$ws3.Cells.Item(10, 4).Value = 2
$ws3.Cells.Item(10, 5).Value = 1
$ws3.Cells.Item(10, 6).Value = 0
$ws3.Cells.Item(10, 7).Value = 12
$ws3.Cells.Item(10, 8).Value = 49

$ws3.Cells.Item(11, 1).Value = "BOONE"
$ws3.Cells.Item(11, 2).Value = 2046
$ws3.Cells.Item(11, 3).Value = 205
$ws3.Cells.Item(11, 4).Value = 9
$ws3.Cells.Item(11, 5).Value = 21
$ws3.Cells.Item(11, 6).Value = 6
$ws3.Cells.Item(11, 7).Value = 403
$ws3.Cells.Item(11, 8).Value = 424

$ws3.Cells.Item(12, 1).Value = "BROWN"
$ws3.Cells.Item(12, 2).Value = 291
$ws3.Cells.Item(12, 3).Value = 5
$ws3.Cells.Item(12, 4).Value = 0
$ws3.Cells.Item(12, 5).Value = 1
$ws3.Cells.Item(12, 6).Value = 0
$ws3.Cells.Item(12, 7).Value = 2
$ws3.Cells.Item(12, 8).Value = 15

$ws3.Cells.Item(13, 1).Value = "BUREAU"
$ws3.Cells.Item(13, 2).Value = 2025
$ws3.Cells.Item(13, 3).Value = 90
$ws3.Cells.Item(13, 4).Value = 12
$ws3.Cells.Item(13, 5).Value = 12
$ws3.Cells.Item(13, 6).Value = 1
$ws3.Cells.Item(13, 7).Value = 142
$ws3.Cells.Item(13, 8).Value = 150

# ----- Age 66+ -----
$ws4.Cells.Item(3, 1).Value = "ALEXANDER"
$ws4.Cells.Item(3, 2).Value = 84
$ws4.Cells.Item(3, 3).Value = 78
$ws4.Cells.Item(3, 4).Value = 0
$ws4.Cells.Item(3, 5).Value = 0
$ws4.Cells.Item(3, 6).Value = 0
$ws4.Cells.Item(3, 7).Value = 0
$ws4.Cells.Item(3, 8).Value = 3

$ws4.Cells.Item(4, 1).Value = "BOND"
$ws4.Cells.Item(4, 2).Value = 130
$ws4.Cells.Item(4, 3).Value = 7
$ws4.Cells.Item(4, 4).Value = 0
$ws4.Cells.Item(4, 5).Value = 0
$ws4.Cells.Item(4, 6).Value = 0
$ws4.Cells.Item(4, 7).Value = 2
$ws4.Cells.Item(4, 8).Value = 1

$ws4.Cells.Item(5, 1).Value = "BOONE"
$ws4.Cells.Item(5, 2).Value = 366
$ws4.Cells.Item(5, 3).Value = 3
$ws4.Cells.Item(5, 4).Value = 0
$ws4.Cells.Item(5, 5).Value = 8
$ws4.Cells.Item(5, 6).Value = 0
$ws4.Cells.Item(5, 7).Value = 127
$ws4.Cells.Item(5, 8).Value = 48

$ws4.Cells.Item(6, 1).Value = "BROWN"
$ws4.Cells.Item(6, 2).Value = 27
$ws4.Cells.Item(6, 3).Value = 0
$ws4.Cells.Item(6, 4).Value = 0
$ws4.Cells.Item(6, 5).Value = 0
$ws4.Cells.Item(6, 6).Value = 0
$ws4.Cells.Item(6, 7).Value = 0
$ws4.Cells.Item(6, 8).Value = 0

$ws4.Cells.Item(7, 1).Value = "BUREAU"
$ws4.Cells.Item(7, 2).Value = 300
$ws4.Cells.Item(7, 3).Value = 2
$ws4.Cells.Item(7, 4).Value = 1
$ws4.Cells.Item(7, 5).Value = 8
$ws4.Cells.Item(7, 6).Value = 0
$ws4.Cells.Item(7, 7).Value = 28
$ws4.Cells.Item(7, 8).Value = 11

Write-Host "edit applied"
